$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row -> (nombre_aides, montant_total)
$updates = @{
    2  = @("150", "330000.00")
    3  = @("802", "2010761.80")
    6  = @("16", "78500.00")
    9  = @("30", "66000.00")
    10 = @("243", "565990.66")
    11 = @("102", "314891.77")
    12 = @("23", "94000.00")
    21 = @("45", "111500.00")
    23 = @("100", "332600.00")
    25 = @("15", "31500.00")
    33 = @("419", "1026811.79")
    72 = @("734", "1892946.83")
    74 = @("86", "335000.00")
    83 = @("77", "173200.00")
    84 = @("338", "814972.09")
    85 = @("133", "414192.00")
    86 = @("38", "141709.01")
    87 = @("7", "27500.00")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    $cellC = $ws.Cells.Item($row, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $vals[0]

    $cellD = $ws.Cells.Item($row, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $vals[1]
}

$wb.Save()
